# The author opened the workbook, widened column A to best-fit its
# contents (the "id_building_type" header / numeric ids), and left the
# cell selection on B9 when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-fit column A to its content (header "id_building_type" / numeric ids).
# ColumnWidth is expressed in characters; 15.83 is the closest achievable
# input that reproduces the saved best-fit width for column A.
$ws.Columns("A").ColumnWidth = 15.83

# Leave the active selection on B9, matching the saved view state.
$ws.Range("B9").Select()
